# Scheduled-runner data refresh: updates market-price / profit figures
# (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H-N) for
# specific leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 191.88235
$ws.Range("I33").Value = 191.88235
$ws.Range("K33").Value = 191.88235
$ws.Range("M33").Value = 37.11765
# Row 51
$ws.Range("H51").Value = 7000
$ws.Range("I51").Value = 6500
$ws.Range("K51").Value = 6500
$ws.Range("M51").Value = -6016
# Row 100
$ws.Range("H100").Value = 2186.4375
$ws.Range("I100").Value = 1628.3
$ws.Range("K100").Value = 1628.3
$ws.Range("M100").Value = -1087.3
# Row 107
$ws.Range("H107").Value = 806.3913
$ws.Range("I107").Value = 795.2857
$ws.Range("J107").Value = 923
$ws.Range("K107").Value = 795.2857
$ws.Range("L107").Value = 923
$ws.Range("M107").Value = 1124.7143
$ws.Range("N107").Value = -4763
# Row 111
$ws.Range("H111").Value = 6962.6665
$ws.Range("J111").Value = 2000
$ws.Range("L111").Value = 6000
$ws.Range("N111").Value = -12134
# Row 121
$ws.Range("H121").Value = 15199.375
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 17256.428
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 51769.284
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -55263.284
# Row 129
$ws.Range("H129").Value = 213482.48
$ws.Range("J129").Value = 257190.52
$ws.Range("L129").Value = 771571.5599999999
$ws.Range("N129").Value = -781571.5599999999
# Row 141
$ws.Range("H141").Value = 2498
$ws.Range("I141").Value = 2236.9
$ws.Range("J141").Value = 3368.3333
$ws.Range("K141").Value = 6710.700000000001
$ws.Range("L141").Value = 10104.9999
$ws.Range("M141").Value = -1530.700000000001
$ws.Range("N141").Value = -20464.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 813.2917
$ws.Range("I2").Value = 868.1111
$ws.Range("J2").Value = 648.8333
$ws.Range("K2").Value = 868.1111
$ws.Range("L2").Value = 648.8333
$ws.Range("M2").Value = -755.1111
$ws.Range("N2").Value = -874.8333
# Row 32
$ws.Range("H32").Value = 6414.338
$ws.Range("I32").Value = 5566.0195
$ws.Range("K32").Value = 5566.0195
$ws.Range("M32").Value = -5279.0195
# Row 110
$ws.Range("H110").Value = 1066.25
$ws.Range("I110").Value = 996.4286
$ws.Range("J110").Value = 1555
$ws.Range("K110").Value = 996.4286
$ws.Range("L110").Value = 1555
$ws.Range("M110").Value = 1048.5714
$ws.Range("N110").Value = -5645
# Row 111
$ws.Range("H111").Value = 25000
$ws.Range("J111").Value = 25000
$ws.Range("L111").Value = 25000
$ws.Range("N111").Value = -33180
# Row 112
$ws.Range("H112").Value = 32019
$ws.Range("J112").Value = 32019
$ws.Range("L112").Value = 32019
$ws.Range("N112").Value = -34973
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
# Row 114
$ws.Range("H114").Value = 30323.25
$ws.Range("J114").Value = 30323.25
$ws.Range("L114").Value = 30323.25
$ws.Range("N114").Value = -39001.25
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""
# Row 116
$ws.Range("H116").Value = 813.2917
$ws.Range("I116").Value = 868.1111
$ws.Range("J116").Value = 648.8333
$ws.Range("K116").Value = 868.1111
$ws.Range("L116").Value = 648.8333
$ws.Range("M116").Value = 1425.8889
$ws.Range("N116").Value = -5236.8333
# Row 117
$ws.Range("H117").Value = 10000
$ws.Range("J117").Value = 10000
$ws.Range("L117").Value = 10000
$ws.Range("N117").Value = -19178
# Row 119
$ws.Range("H119").Value = 30250
$ws.Range("J119").Value = 30250
$ws.Range("L119").Value = 30250
$ws.Range("N119").Value = -39926
# Row 121
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
# Row 131
$ws.Range("H131").Value = 60049.57
$ws.Range("J131").Value = 60049.57
$ws.Range("L131").Value = 60049.57
$ws.Range("N131").Value = -70129.57000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 813.2917
$ws.Range("I3").Value = 868.1111
$ws.Range("J3").Value = 648.8333
$ws.Range("K3").Value = 868.1111
$ws.Range("L3").Value = 648.8333
$ws.Range("M3").Value = -754.1111
$ws.Range("N3").Value = -876.8333
# Row 99
$ws.Range("H99").Value = 1622.1111
$ws.Range("I99").Value = 1728.5714
$ws.Range("K99").Value = 1728.5714
$ws.Range("M99").Value = -230.5714

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 28998
$ws.Range("J52").Value = 28998
$ws.Range("L52").Value = 28998
$ws.Range("N52").Value = -29586
# Row 58
$ws.Range("H58").Value = 17459.936
$ws.Range("I58").Value = 1289.6154
$ws.Range("K58").Value = 1289.6154
$ws.Range("M58").Value = -1086.6154
# Row 110
$ws.Range("H110").Value = 30702
$ws.Range("J110").Value = 30702
$ws.Range("L110").Value = 30702
$ws.Range("N110").Value = -38882
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""
# Row 136
$ws.Range("H136").Value = 17459.936
$ws.Range("I136").Value = 1289.6154
$ws.Range("K136").Value = 3868.8462
$ws.Range("M136").Value = -1318.8462

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 729.09
$ws.Range("J131").Value = 730.0909
$ws.Range("L131").Value = 2190.2727
$ws.Range("N131").Value = -12270.2727
# Row 141
$ws.Range("H141").Value = 5488.3335
$ws.Range("J141").Value = 5488.3335
$ws.Range("L141").Value = 16465.0005
$ws.Range("N141").Value = -26825.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2627.558
$ws.Range("I126").Value = 2198.6
$ws.Range("J126").Value = 3000.5652
$ws.Range("K126").Value = 6595.799999999999
$ws.Range("L126").Value = 9001.695599999999
$ws.Range("M126").Value = -4125.799999999999
$ws.Range("N126").Value = -13941.6956

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 717.1429000000001
$ws.Range("I93").Value = 717.1429000000001
$ws.Range("K93").Value = 717.1429000000001
$ws.Range("M93").Value = 530.8570999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 142858160
$ws.Range("I107").Value = 250000240
$ws.Range("J107").Value = 2033.3334
$ws.Range("K107").Value = 750000720
$ws.Range("L107").Value = 6100.0002
$ws.Range("M107").Value = -749998800
$ws.Range("N107").Value = -9940.0002
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
# Row 112
$ws.Range("H112").Value = 21250
$ws.Range("J112").Value = 21250
$ws.Range("L112").Value = 21250
$ws.Range("N112").Value = -24204
# Row 113
$ws.Range("H113").Value = 870.2857
$ws.Range("I113").Value = 883.3158
$ws.Range("J113").Value = 746.5
$ws.Range("K113").Value = 2649.9474
$ws.Range("L113").Value = 2239.5
$ws.Range("M113").Value = -479.9474
$ws.Range("N113").Value = -6579.5
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
